$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-24 Tuesday" "2024-09-25 Wednesday"

Replace-Text "595÷9=66, 1" "217÷3=72, 1"
Replace-Text "679÷5=135, 4" "762÷2=381, 0"
Replace-Text "336÷2=168, 0" "207÷7=29, 4"
Replace-Text "997÷2=498, 1" "405÷7=57, 6"
Replace-Text "618÷4=154, 2" "776÷9=86, 2"

Replace-Text "662÷5=132, 2" "862÷8=107, 6"
Replace-Text "331÷4=82, 3" "461÷4=115, 1"
Replace-Text "586÷6=97, 4" "437÷2=218, 1"
Replace-Text "439÷6=73, 1" "335÷8=41, 7"
Replace-Text "869÷7=124, 1" "310÷9=34, 4"

Replace-Text "457÷2=228, 1" "735÷9=81, 6"
Replace-Text "428÷8=53, 4" "598÷3=199, 1"
Replace-Text "143÷8=17, 7" "598÷9=66, 4"
Replace-Text "409÷5=81, 4" "218÷5=43, 3"
Replace-Text "595÷6=99, 1" "843÷8=105, 3"

Replace-Text "794÷3=264, 2" "296÷6=49, 2"
Replace-Text "671÷8=83, 7" "178÷3=59, 1"
Replace-Text "591÷3=197, 0" "810÷5=162, 0"
Replace-Text "210÷6=35, 0" "144÷2=72, 0"
Replace-Text "371÷5=74, 1" "437÷5=87, 2"

Replace-Text "957÷9=106, 3" "873÷3=291, 0"
Replace-Text "244÷2=122, 0" "345÷9=38, 3"
Replace-Text "261÷2=130, 1" "967÷9=107, 4"
Replace-Text "126÷8=15, 6" "441÷3=147, 0"
Replace-Text "666÷5=133, 1" "180÷7=25, 5"
